# Apply Betfair Back/Lay odds updates for 2025-10-24 workbook.
# Each assignment below corresponds to one changed cell in the source diff
# (old value -> new value shown in the trailing comment).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value2 = 3.15  # G2: 3.2 -> 3.15
$ws.Cells.Item(2, 8).Value2 = 2.26  # H2: 2.24 -> 2.26
$ws.Cells.Item(2, 9).Value2 = 2.36  # I2: 2.38 -> 2.36
$ws.Cells.Item(2, 10).Value2 = 4  # J2: 3.9 -> 4
$ws.Cells.Item(2, 12).Value2 = 1.28  # L2: 1.23 -> 1.28
$ws.Cells.Item(2, 14).Value2 = 5.9  # N2: 5.6 -> 5.9
$ws.Cells.Item(2, 16).Value2 = 2.7  # P2: 2.58 -> 2.7
$ws.Cells.Item(2, 17).Value2 = 1.53  # Q2: 1.55 -> 1.53
$ws.Cells.Item(2, 18).Value2 = 1.7  # R2: 1.64 -> 1.7
$ws.Cells.Item(2, 19).Value2 = 2.3  # S2: 2.34 -> 2.3
$ws.Cells.Item(2, 20).Value2 = 1.51  # T2: 1.53 -> 1.51
$ws.Cells.Item(2, 21).Value2 = 2.72  # U2: 2.62 -> 2.72
$ws.Cells.Item(2, 23).Value2 = 1.46  # W2: 1.45 -> 1.46
$ws.Cells.Item(2, 24).Value2 = 1000  # X2: 32 -> 1000
$ws.Cells.Item(2, 25).Value2 = 1000  # Y2: 19 -> 1000
$ws.Cells.Item(2, 26).Value2 = 1000  # Z2: 23 -> 1000
$ws.Cells.Item(2, 27).Value2 = 1000  # AA2: 36 -> 1000
$ws.Cells.Item(2, 28).Value2 = 1000  # AB2: 23 -> 1000
$ws.Cells.Item(2, 29).Value2 = 10.5  # AC2: 12.5 -> 10.5
$ws.Cells.Item(2, 30).Value2 = 12.5  # AD2: 14 -> 12.5
$ws.Cells.Item(2, 31).Value2 = 1000  # AE2: 25 -> 1000
$ws.Cells.Item(2, 32).Value2 = 1000  # AF2: 30 -> 1000
$ws.Cells.Item(2, 33).Value2 = 14.5  # AG2: 17 -> 14.5
$ws.Cells.Item(2, 34).Value2 = 15  # AH2: 18 -> 15
$ws.Cells.Item(2, 35).Value2 = 1000  # AI2: 32 -> 1000
$ws.Cells.Item(2, 36).Value2 = 1000  # AJ2: 60 -> 1000
$ws.Cells.Item(2, 37).Value2 = 1000  # AK2: 34 -> 1000
$ws.Cells.Item(2, 38).Value2 = 1000  # AL2: 40 -> 1000
$ws.Cells.Item(2, 39).Value2 = 1000  # AM2: 65 -> 1000
$ws.Cells.Item(2, 40).Value2 = 1000  # AN2: 21 -> 1000
$ws.Cells.Item(2, 41).Value2 = 11  # AO2: 13.5 -> 11

# Row 3
$ws.Cells.Item(3, 10).Value2 = 3.55  # J3: 3.6 -> 3.55
$ws.Cells.Item(3, 16).Value2 = 2.12  # P3: 2.14 -> 2.12

# Row 4
$ws.Cells.Item(4, 10).Value2 = 3.65  # J4: 3.7 -> 3.65
$ws.Cells.Item(4, 16).Value2 = 2.08  # P4: 2.24 -> 2.08
$ws.Cells.Item(4, 17).Value2 = 1.62  # Q4: 1.63 -> 1.62

# Row 5
$ws.Cells.Item(5, 16).Value2 = 2.06  # P5: 2.1 -> 2.06
$ws.Cells.Item(5, 17).Value2 = 1.64  # Q5: 1.77 -> 1.64

# Row 6
$ws.Cells.Item(6, 6).Value2 = 3.9  # F6: 3.7 -> 3.9
$ws.Cells.Item(6, 7).Value2 = 5.2  # G6: 4.4 -> 5.2
$ws.Cells.Item(6, 8).Value2 = 1.72  # H6: 1.79 -> 1.72
$ws.Cells.Item(6, 9).Value2 = 1.89  # I6: 1.94 -> 1.89
$ws.Cells.Item(6, 10).Value2 = 4.4  # J6: 3.95 -> 4.4
$ws.Cells.Item(6, 11).Value2 = 5.3  # K6: 5.2 -> 5.3
$ws.Cells.Item(6, 16).Value2 = 3.15  # P6: 3.1 -> 3.15
$ws.Cells.Item(6, 19).Value2 = 1.8  # S6: 1.92 -> 1.8
$ws.Cells.Item(6, 20).Value2 = 1.45  # T6: 1.44 -> 1.45
$ws.Cells.Item(6, 21).Value2 = 2.52  # U6: 2.56 -> 2.52
$ws.Cells.Item(6, 29).Value2 = 1000  # AC6: 15.5 -> 1000
$ws.Cells.Item(6, 30).Value2 = 1000  # AD6: 14 -> 1000
$ws.Cells.Item(6, 41).Value2 = 6.2  # AO6: 7.4 -> 6.2

# Row 7
$ws.Cells.Item(7, 19).Value2 = 3  # S7: 3.05 -> 3
$ws.Cells.Item(7, 20).Value2 = 1.68  # T7: 1.69 -> 1.68
$ws.Cells.Item(7, 21).Value2 = 2.36  # U7: 2.32 -> 2.36
$ws.Cells.Item(7, 25).Value2 = 15.5  # Y7: 17 -> 15.5
$ws.Cells.Item(7, 26).Value2 = 27  # Z7: 30 -> 27
$ws.Cells.Item(7, 28).Value2 = 11.5  # AB7: 12.5 -> 11.5
$ws.Cells.Item(7, 30).Value2 = 15.5  # AD7: 16.5 -> 15.5
$ws.Cells.Item(7, 31).Value2 = 46  # AE7: 40 -> 46
$ws.Cells.Item(7, 34).Value2 = 16.5  # AH7: 18.5 -> 16.5
$ws.Cells.Item(7, 35).Value2 = 46  # AI7: 50 -> 46
$ws.Cells.Item(7, 36).Value2 = 30  # AJ7: 32 -> 30
$ws.Cells.Item(7, 37).Value2 = 24  # AK7: 23 -> 24
$ws.Cells.Item(7, 39).Value2 = 1000  # AM7: 95 -> 1000
$ws.Cells.Item(7, 40).Value2 = 14.5  # AN7: 15 -> 14.5
$ws.Cells.Item(7, 41).Value2 = 34  # AO7: 36 -> 34

# Row 8
$ws.Cells.Item(8, 6).Value2 = 1.75  # F8: 1.73 -> 1.75
$ws.Cells.Item(8, 8).Value2 = 5.3  # H8: 5.4 -> 5.3
$ws.Cells.Item(8, 9).Value2 = 5.5  # I8: 5.6 -> 5.5
$ws.Cells.Item(8, 24).Value2 = 16.5  # X8: 21 -> 16.5
$ws.Cells.Item(8, 25).Value2 = 20  # Y8: 22 -> 20
$ws.Cells.Item(8, 26).Value2 = 55  # Z8: 46 -> 55
$ws.Cells.Item(8, 29).Value2 = 9.4  # AC8: 9.199999999999999 -> 9.4
$ws.Cells.Item(8, 30).Value2 = 22  # AD8: 25 -> 22
$ws.Cells.Item(8, 31).Value2 = 1000  # AE8: 85 -> 1000
$ws.Cells.Item(8, 36).Value2 = 18  # AJ8: 20 -> 18
$ws.Cells.Item(8, 37).Value2 = 18  # AK8: 20 -> 18
$ws.Cells.Item(8, 38).Value2 = 36  # AL8: 42 -> 36

# Row 9
$ws.Cells.Item(9, 6).Value2 = 1.36  # F9: 1.34 -> 1.36
$ws.Cells.Item(9, 7).Value2 = 1.53  # G9: 1.55 -> 1.53
$ws.Cells.Item(9, 8).Value2 = 2.86  # H9: 2.78 -> 2.86
$ws.Cells.Item(9, 10).Value2 = 2.86  # J9: 2.78 -> 2.86

# Row 10
$ws.Cells.Item(10, 6).Value2 = 1.92  # F10: 1.86 -> 1.92
$ws.Cells.Item(10, 7).Value2 = 2.32  # G10: 2.34 -> 2.32
$ws.Cells.Item(10, 8).Value2 = 3.55  # H10: 3.5 -> 3.55
$ws.Cells.Item(10, 9).Value2 = 5  # I10: 5.1 -> 5
$ws.Cells.Item(10, 10).Value2 = 2.68  # J10: 2.64 -> 2.68
$ws.Cells.Item(10, 11).Value2 = 5.2  # K10: 5.8 -> 5.2

# Row 11
$ws.Cells.Item(11, 7).Value2 = 1.32  # G11: 1.33 -> 1.32
$ws.Cells.Item(11, 8).Value2 = 11.5  # H11: 12 -> 11.5
$ws.Cells.Item(11, 9).Value2 = 13.5  # I11: 14 -> 13.5
$ws.Cells.Item(11, 10).Value2 = 6.2  # J11: 6 -> 6.2
$ws.Cells.Item(11, 14).Value2 = 4.7  # N11: 4.8 -> 4.7
$ws.Cells.Item(11, 18).Value2 = 1.48  # R11: 1.5 -> 1.48
$ws.Cells.Item(11, 19).Value2 = 2.86  # S11: 2.8 -> 2.86
$ws.Cells.Item(11, 24).Value2 = 22  # X11: 21 -> 22
$ws.Cells.Item(11, 34).Value2 = 38  # AH11: 36 -> 38
$ws.Cells.Item(11, 36).Value2 = 10  # AJ11: 10.5 -> 10
$ws.Cells.Item(11, 37).Value2 = 15  # AK11: 15.5 -> 15
$ws.Cells.Item(11, 39).Value2 = 250  # AM11: 240 -> 250
$ws.Cells.Item(11, 40).Value2 = 5.3  # AN11: 5.2 -> 5.3

# Row 12
$ws.Cells.Item(12, 6).Value2 = 1.91  # F12: 1.9 -> 1.91
$ws.Cells.Item(12, 7).Value2 = 1.92  # G12: 1.91 -> 1.92
$ws.Cells.Item(12, 14).Value2 = 3.6  # N12: 3.65 -> 3.6
$ws.Cells.Item(12, 16).Value2 = 1.89  # P12: 1.88 -> 1.89
$ws.Cells.Item(12, 18).Value2 = 1.34  # R12: 1.33 -> 1.34
$ws.Cells.Item(12, 20).Value2 = 1.93  # T12: 1.92 -> 1.93
$ws.Cells.Item(12, 21).Value2 = 1.99  # U12: 2 -> 1.99
$ws.Cells.Item(12, 27).Value2 = 130  # AA12: 120 -> 130
$ws.Cells.Item(12, 28).Value2 = 8.6  # AB12: 8.4 -> 8.6
$ws.Cells.Item(12, 32).Value2 = 10.5  # AF12: 11 -> 10.5
$ws.Cells.Item(12, 35).Value2 = 80  # AI12: 75 -> 80
$ws.Cells.Item(12, 36).Value2 = 21  # AJ12: 20 -> 21
$ws.Cells.Item(12, 38).Value2 = 38  # AL12: 42 -> 38
$ws.Cells.Item(12, 40).Value2 = 14  # AN12: 14.5 -> 14
$ws.Cells.Item(12, 41).Value2 = 85  # AO12: 90 -> 85

# Row 13
$ws.Cells.Item(13, 6).Value2 = 2.1  # F13: 2.08 -> 2.1
$ws.Cells.Item(13, 7).Value2 = 2.14  # G13: 2.12 -> 2.14
$ws.Cells.Item(13, 8).Value2 = 4  # H13: 4.1 -> 4
$ws.Cells.Item(13, 9).Value2 = 4.3  # I13: 4.4 -> 4.3
$ws.Cells.Item(13, 11).Value2 = 3.55  # K13: 3.5 -> 3.55
$ws.Cells.Item(13, 14).Value2 = 3.5  # N13: 3.55 -> 3.5
$ws.Cells.Item(13, 18).Value2 = 1.32  # R13: 1.33 -> 1.32
$ws.Cells.Item(13, 19).Value2 = 3.9  # S13: 3.85 -> 3.9
$ws.Cells.Item(13, 27).Value2 = 95  # AA13: 110 -> 95
$ws.Cells.Item(13, 32).Value2 = 12.5  # AF13: 13 -> 12.5
$ws.Cells.Item(13, 34).Value2 = 20  # AH13: 19.5 -> 20
$ws.Cells.Item(13, 37).Value2 = 24  # AK13: 23 -> 24
$ws.Cells.Item(13, 38).Value2 = 44  # AL13: 42 -> 44
$ws.Cells.Item(13, 39).Value2 = 140  # AM13: 110 -> 140
$ws.Cells.Item(13, 40).Value2 = 18  # AN13: 17.5 -> 18
$ws.Cells.Item(13, 41).Value2 = 80  # AO13: 70 -> 80

# Row 14
$ws.Cells.Item(14, 7).Value2 = 2.04  # G14: 2.08 -> 2.04
$ws.Cells.Item(14, 8).Value2 = 4.9  # H14: 4.8 -> 4.9

# Row 15
$ws.Cells.Item(15, 7).Value2 = 1.6  # G15: 1.59 -> 1.6
$ws.Cells.Item(15, 9).Value2 = 10.5  # I15: 9.4 -> 10.5
$ws.Cells.Item(15, 10).Value2 = 4  # J15: 3.75 -> 4
